$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'compression shorts for woman'
    'compression shorts for women'
    'compression shorts for women running'
    'compression shorts hamstring'
    'compression shorts knee length'
    'compression shorts ladies'
    'compression shorts long'
    'compression shorts long women'
    'compression shorts low rise'
    'compression shorts quad'
    'compression shorts with padding basketball'
    'compression shorts women'
    'compression shorts women 2x'
    'compression shorts women 3 inch'
    'compression shorts women 3 pack'
    'compression shorts women 3x'
    'compression shorts women 4 inch'
    'compression shorts women 7 inch'
    'compression shorts women black'
    'compression shorts women high waist'
    'compression shorts women knee length'
    'compression shorts women long'
    'compression shorts women pack'
    'compression shorts women plus'
    'compression shorts women plus size'
    'compression shorts women running'
    'compression shorts women soccer'
    'compression shorts women xl'
    'compression shorts womens'
    'compression shorts womens 3 pack'
    'compression slip'
    'compression spandex'
    'compression sport'
    'compression sport tights'
    'compression sports'
    'compression sports pants'
    'compression stocks'
    'compression stocks women'
    'compression stomach men'
    'compression support'
    'compression support garment'
    'compression support hose'
    'compression support knee'
    'compression support knee highs'
    'compression support shorts'
    'compression tank boys'
    'compression tank for women'
    'compression tank top women long'
    'compression tank tops for women tummy control'
    'compression tank women'
    'compression tanks for women'
    'compression tanks women'
    'compression tape'
    'compression therapy leg'
    'compression therapy pants'
    'compression thermal'
    'compression thermal for men'
    'compression thigh brace'
    'compression thigh support'
    'compression thights'
    'compression tight'
    'compression tight leggings'
    'compression tight men'
    'compression tight pants'
    'compression tight women'
    'compression tights'
    'compression tights 2xu'
    'compression tights 3 4'
    'compression tights blue'
    'compression tights for boys'
    'compression tights for boys basketball'
    'compression tights for flying'
    'compression tights for men 3 pack'
    'compression tights for men cold weather'
    'compression tights for women'
    'compression tights for women circulation'
    'compression tights for women high waist'
    'compression tights medical'
    'compression tights men'
    'compression tights mens'
    'compression tights shorts'
    'compression tights shorts women'
    'compression tights with knee pads'
    'compression tights with pockets'
    'compression tights women'
    'compression tights women 20-30mmhg'
    'compression tights women black'
    'compression tights women circulation'
    'compression tights women cwx'
    'compression tights women medical'
    'compression tights women running'
    'compression tights women white'
    'compression top for women'
    'compression top thermal'
    'compression top women'
    'compression top women workout'
    'compression tops for women'
    'compression tops women'
    'compression travel'
    'compression trouser'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

